## Model Spawn Added. Marker spawn Added.
## Updates the "Main" sheet's test-file path + world file, and the
## "Sheet1" sheet's spawned-model filenames, then restores the on-screen
## selection/column sizing that the original author left behind.

$wb = $excel.ActiveWorkbook

# ---- "Main" sheet -------------------------------------------------------
$wsMain = $wb.Worksheets.Item(1)
$wsMain.Activate() | Out-Null

# test_files_path now points at the Simulations sub-folder.
$wsMain.Range("B2").Value = "/home/stb21753492/FiducialTags/Simulations"

# world_file switched from the single-marker demo world to the standard one.
$wsMain.Range("B3").Value = "standard_world.sdf"

# Widen column B so the longer path/file strings are readable (target stored
# width 38.9 chars), and leave the selection where the author left it after
# editing.
$wsMain.Columns.Item(2).ColumnWidth = 38.0667
$wsMain.Range("B5").Select() | Out-Null

# ---- "Sheet1" sheet -------------------------------------------------------
$wsModels = $wb.Worksheets.Item(2)
$wsModels.Activate() | Out-Null

# Spawned models renamed from the generic placeholders to the real fiducial
# marker model files.
$wsModels.Range("B7").Value = "DICT_4X4_50_s500_id0.sdf"
$wsModels.Range("B8").Value = "DICT_4X4_50_s500_id1.sdf"

$wsModels.Range("B8").Select() | Out-Null
